$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 231, pushing the existing rows
# (231-307) down to (233-309). This mirrors the upstream edit where
# two new weekly price records were added and the rest of the table
# shifted down.
$ws.Range("A231:A232").EntireRow.Insert()

# ---- New row 231 : Acelga, "Primera" quality ----
$ws.Cells.Item(231, 1).Value = 7
$ws.Cells.Item(231, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(231, 3).Value = "Ñuble"
$ws.Cells.Item(231, 4).Value = 44876
$ws.Cells.Item(231, 5).Value = 16
$ws.Cells.Item(231, 6).Value = 100112009
$ws.Cells.Item(231, 7).Value = "Acelga"
$ws.Cells.Item(231, 8).Value = "Sin especificar"
$ws.Cells.Item(231, 9).Value = "Primera"
$ws.Cells.Item(231, 10).Value = 300
$ws.Cells.Item(231, 11).Value = 600
$ws.Cells.Item(231, 12).Value = 700
$ws.Cells.Item(231, 13).Value = 650
$ws.Cells.Item(231, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(231, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(231, 16).Value = 650
$ws.Cells.Item(231, 17).Value = 1
$ws.Cells.Item(231, 18).Value = "Hortaliza"

# ---- New row 232 : Acelga, "Segunda" quality ----
$ws.Cells.Item(232, 1).Value = 7
$ws.Cells.Item(232, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(232, 3).Value = "Ñuble"
$ws.Cells.Item(232, 4).Value = 44876
$ws.Cells.Item(232, 5).Value = 16
$ws.Cells.Item(232, 6).Value = 100112009
$ws.Cells.Item(232, 7).Value = "Acelga"
$ws.Cells.Item(232, 8).Value = "Sin especificar"
$ws.Cells.Item(232, 9).Value = "Segunda"
$ws.Cells.Item(232, 10).Value = 200
$ws.Cells.Item(232, 11).Value = 500
$ws.Cells.Item(232, 12).Value = 500
$ws.Cells.Item(232, 13).Value = 500
$ws.Cells.Item(232, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(232, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(232, 16).Value = 500
$ws.Cells.Item(232, 17).Value = 1
$ws.Cells.Item(232, 18).Value = "Hortaliza"

# Make sure column D keeps the date number format for the two new rows
$ws.Range("D231:D232").NumberFormat = $ws.Range("D233").NumberFormat
